# "cambio de primera clasificacion"
#
# The sheet's header row goes from 12 columns (A1:L1) down to 7 (A1:G1):
#   old: REFERENCIA, PEDIMENTO, CLIENTE, FRACCION, DESCRIPCION, INGLES,
#        MODELO, F DE PREVIO, FACTURA, AUTOR, CONSECUTIVO, REGIMEN
#   new: REFERENCIA, FRACCION, DESCRIPCION, INGLES, MODELO, F DE PREVIO,
#        FACTURA
#
# PEDIMENTO, CLIENTE, AUTOR, CONSECUTIVO and REGIMEN are dropped; the
# survivors are reordered and the autofilter / filter-database range
# shrink from A1:L1 to A1:G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the autofilter before touching the grid so it doesn't keep
# pointing at stale columns while we rework them.
$ws.AutoFilterMode = $false

# FRACCION/DESCRIPCION/INGLES/MODELO/F DE PREVIO/FACTURA stay; shift them
# left onto B:G; then remove the now-empty trailing columns (old H:L --
# F DE PREVIO/FACTURA/AUTOR/CONSECUTIVO/REGIMEN) entirely.
$ws.Range("H1:L1").EntireColumn.Delete()

$ws.Range("A1").Value = "REFERENCIA"
$ws.Range("B1").Value = "FRACCION"
$ws.Range("C1").Value = "DESCRIPCION"
$ws.Range("D1").Value = "INGLES"
$ws.Range("E1").Value = "MODELO"
$ws.Range("F1").Value = "F DE PREVIO"
$ws.Range("G1").Value = "FACTURA"

# Re-apply the autofilter over the new, narrower header.
[void]$ws.Range("A1:G1").AutoFilter()

# The hidden _FilterDatabase defined name tracks the autofilter range but
# isn't refreshed automatically - point it at the new extent explicitly.
$wb.Names.Item(1).RefersTo = "=Hoja1!`$A`$1:`$G`$1"

# Matches the cached selection left behind in the edited workbook.
[void]$ws.Range("F10").Select()
